$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# K1: date changed from 2018-02-09 (43140) to 2018-03-06 (43165)
$ws.Range("K1").Value = 43165

# Widen column A (was 31.42578125 chars ~ add 5 chars -> 36.42578125)
$ws.Columns.Item(1).ColumnWidth = 35.6

# New task row 16: "PHP Code SQL injection sicher machen", marked "S" in column K
$ws.Range("A16").Value = "PHP Code SQL injection sicher machen"
$ws.Range("K16").Value = "S"

# Update the remembered selection/active cell shown in the sheet view
$ws.Range("I21").Select()
